$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 0
    3 = 0
    4 = 2
    5 = 1
    6 = 1
    7 = 1
    8 = 0
    9 = 0
    10 = 1
    11 = 3
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 2
    17 = 3
    18 = 1
    19 = 1
    20 = 0
    21 = 3
    22 = 1
    23 = 0
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 1
    29 = 2
    30 = 1
    31 = 2
    32 = 1
    33 = 2
    34 = 1
    36 = 1
    37 = 1
    38 = 1
    39 = 0
    40 = 1
    41 = 2
    42 = 0
    43 = 2
    44 = 0
    45 = 2
    46 = 0
    47 = 1
    48 = 0
    49 = 3
    50 = 1
    51 = 1
    52 = 2
    53 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
